$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: wrap row-11 "setText" formulas with a zero-check guard
$ws.Range("F16").Formula = '="if(board["&LEFT(A1,1)-1&"]["&RIGHT(A1,1)-1&"]!=0) "&F11'
$ws.Range("G16:I16").Formula = '="if(board["&LEFT(B1,1)-1&"]["&RIGHT(B1,1)-1&"]!=0) "&G11'

# Row 17: wrap row-12 "setText" formulas with a zero-check guard
$ws.Range("F17:F19").Formula = '="if(board["&LEFT(A2,1)-1&"]["&RIGHT(A2,1)-1&"]!=0) "&F12'
$ws.Range("G17:G19").Formula = '="if(board["&LEFT(B2,1)-1&"]["&RIGHT(B2,1)-1&"]!=0) "&G12'
$ws.Range("H17:H19").Formula = '="if(board["&LEFT(C2,1)-1&"]["&RIGHT(C2,1)-1&"]!=0) "&H12'
$ws.Range("I17:I19").Formula = '="if(board["&LEFT(D2,1)-1&"]["&RIGHT(D2,1)-1&"]!=0) "&I12'

$ws.Range("F16:I19").Select() | Out-Null
